$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.022.03"
$ws.Range("E2").Value = "  -2.97%  "
$ws.Range("D3").Value = "1.726.97"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "309.96"
$ws.Range("E5").Value = "  -5.41%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "0.4807"
$ws.Range("E7").Value = "  +3.64%  "
$ws.Range("D8").Value = "0.3480"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("D9").Value = "43.34"
$ws.Range("E9").Value = "  +3.28%  "
$ws.Range("D10").Value = "0.07236"
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("D11").Value = "1.050"
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "19.92"
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").Value = "5.871"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "1.728.46"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "6.857"
$ws.Range("E16").Value = "  -3.96%  "
$ws.Range("D17").Value = "86.82"
$ws.Range("E17").Value = "  -5.72%  "
$ws.Range("D18").Value = "0.00001034"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").Value = "0.06382"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "16.65"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").Value = "5.713"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").Value = "27.076.29"
$ws.Range("E23").Value = "  -2.81%  "
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").Value = "2.084"
$ws.Range("E25").Value = "  -3.12%  "
$ws.Range("D26").Value = "154.40"
$ws.Range("E26").Value = "  -4.44%  "
$ws.Range("D27").Value = "20.00"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "1.925.96"
$ws.Range("D29").Value = "2.066"
$ws.Range("E29").Value = "  -3.96%  "
$ws.Range("D30").Value = "120.67"
$ws.Range("E30").Value = "  -1.59%  "
$ws.Range("D31").Value = "1.045"
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("D32").Value = "0.09308"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "3.648"
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").Value = "5.374"
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("D35").Value = "0.05940"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("D36").Value = "0.02185"
$ws.Range("E36").Value = "  -3.48%  "
$ws.Range("B37").Value = "WEMIXTOKEN"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "1.429"
$ws.Range("E37").Value = "  +6.20%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "10.97"
$ws.Range("E38").Value = "  -5.52%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "4.765"
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "0.1992"
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").Value = "0.5975"
$ws.Range("E42").Value = "  -2.97%  "
$ws.Range("D43").Value = "1.095"
$ws.Range("E43").Value = "  -6.87%  "
$ws.Range("D44").Value = "7.536"
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("D45").Value = "12.67"
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("D47").Value = "0.5619"
$ws.Range("E47").Value = "  -2.82%  "
$ws.Range("D48").Value = "118.68"
$ws.Range("E48").Value = "  -3.21%  "
$ws.Range("D49").Value = "1.847"
$ws.Range("E49").Value = "  -3.90%  "
$ws.Range("D50").Value = "1.100"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").Value = "0.06642"
$ws.Range("E51").Value = "  -2.21%  "
